$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$dates = @(
    "Tue Jul 11 12:57:42 EDT 2023",
    "Tue Jul 11 12:57:50 EDT 2023",
    "Tue Jul 11 12:57:58 EDT 2023",
    "Tue Jul 11 12:58:07 EDT 2023",
    "Tue Jul 11 12:58:16 EDT 2023",
    "Tue Jul 11 12:58:25 EDT 2023",
    "Tue Jul 11 12:58:33 EDT 2023",
    "Tue Jul 11 12:58:42 EDT 2023",
    "Tue Jul 11 12:58:51 EDT 2023",
    "Tue Jul 11 12:58:59 EDT 2023",
    "Tue Jul 11 12:59:08 EDT 2023",
    "Tue Jul 11 12:59:17 EDT 2023",
    "Tue Jul 11 12:59:26 EDT 2023",
    "Tue Jul 11 12:59:34 EDT 2023",
    "Tue Jul 11 12:59:43 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $dates[$i]
}
